# Append 9 new NBA game rows (763-771) to the bottom of the existing table
# on sheet1, matching columns:
#   A: Away team   B: Away Pts   C: Home team   D: Home Pts
#   E: Overtime    F: Attend.    G: Arena        H: Win team
#   I: Loss team

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @("Golden State Warriors", 131, "Indiana Pacers", 109, "No", 17832, "Gainbridge Fieldhouse", "Golden State Warriors", "Indiana Pacers"),
    @("San Antonio Spurs", 111, "Orlando Magic", 127, "No", 17832, "Amway Center", "Orlando Magic", "San Antonio Spurs"),
    @("Cleveland Cavaliers", 118, "Brooklyn Nets", 95, "No", 17832, "Barclays Center", "Cleveland Cavaliers", "Brooklyn Nets"),
    @("Dallas Mavericks", 122, "New York Knicks", 108, "No", 17832, "Madison Square Garden (IV)", "Dallas Mavericks", "New York Knicks"),
    @("Chicago Bulls", 118, "Memphis Grizzlies", 110, "No", 17832, "FedEx Forum", "Chicago Bulls", "Memphis Grizzlies"),
    @("Minnesota Timberwolves", 129, "Milwaukee Bucks", 105, "No", 17832, "Fiserv Forum", "Minnesota Timberwolves", "Milwaukee Bucks"),
    @("Utah Jazz", 115, "Phoenix Suns", 129, "No", 17832, "Footprint Center", "Phoenix Suns", "Utah Jazz"),
    @("Denver Nuggets", 114, "Los Angeles Lakers", 106, "No", 17832, "Crypto.com Arena", "Denver Nuggets", "Los Angeles Lakers"),
    @("Detroit Pistons", 128, "Portland Trail Blazers", 122, "OT", 17832, "Moda Center", "Detroit Pistons", "Portland Trail Blazers")
)

$startRow = 763
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
    # Match the existing number format used by the Pts columns (B, D)
    $ws.Cells.Item($r, 2).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 4).NumberFormat = "#,##0"
}

# Update the view to match the post-edit scroll/selection position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 739
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G" + $startRow).Select()
